$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MIBGAS")

$ws.Range("B2:B419").Value = 30.08
